# Add notification mailing system entries to the error/message list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 239; Num = 35; Category = "Notify"; Message = "NOTIFICATION" },
    @{ Row = 240; Num = 34; Category = "Notify"; Message = "Past Pick-ups" },
    @{ Row = 241; Num = 33; Category = "Notify"; Message = "Upcoming Trees" },
    @{ Row = 242; Num = 32; Category = "Notify"; Message = "Notify Checked Trees to Managers" },
    @{ Row = 243; Num = 31; Category = "Notify"; Message = "Notify Checked Trees to Parents" },
    @{ Row = 244; Num = 30; Category = "Notify"; Message = "Notification has sent successfully." },
    @{ Row = 245; Num = 29; Category = "Notify"; Message = "Failed to send a notification to parents." },
    @{ Row = 246; Num = 28; Category = "Notify"; Message = "Failed to send a notification to managers." },
    @{ Row = 247; Num = 27; Category = "Notify"; Message = "Server internal error. Please try again." }
)

foreach ($r in $rows) {
    $ws.Range("A" + $r.Row).Value = $r.Num
    $ws.Range("B" + $r.Row).Value = $r.Category
    $ws.Range("C" + $r.Row).Value = $r.Message
    $ws.Range("D" + $r.Row).Value = "Message"
}

$ws.Range("C244").Select() | Out-Null
